# Rotate the species-observation data among rows 11, 12 and 13:
#   new row 11 <- old row 12
#   new row 12 <- old row 13
#   new row 13 <- old row 11
# Only columns A, B, E, F, G, H, Q, R, Z, AB are affected; every other
# column is identical across these three rows so it is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the "old" values for rows 11-13 before overwriting anything.
$cols = @("A","B","E","F","G","H","Q","R","Z","AB")

$old11 = @{}
$old12 = @{}
$old13 = @{}
foreach ($col in $cols) {
    $old11[$col] = $ws.Range("$col`11").Value2
    $old12[$col] = $ws.Range("$col`12").Value2
    $old13[$col] = $ws.Range("$col`13").Value2
}

foreach ($col in $cols) {
    $v = $old12[$col]
    if ($null -eq $v) { $v = "" }
    $ws.Range("$col`11").Value = $v
}
foreach ($col in $cols) {
    $v = $old13[$col]
    if ($null -eq $v) { $v = "" }
    $ws.Range("$col`12").Value = $v
}
foreach ($col in $cols) {
    $v = $old11[$col]
    if ($null -eq $v) { $v = "" }
    $ws.Range("$col`13").Value = $v
}
